$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last existing row (11) currently carries the special "last row" date
# format (YYYY-MM-DD). Since we are appending a new last row, that special
# format moves down to the new row, and row 11 reverts to the regular
# date/time format used by the other data rows.
$ws.Range("A11").NumberFormat = $ws.Range("A10").NumberFormat

# Append the new data row (2021-11-23).
$ws.Range("A12").Value = 44523
$ws.Range("B12").Value = -98.04999999999973

# The newly appended row becomes the new "last row" and takes on the
# special date-only format previously used by row 11.
$ws.Range("A12").NumberFormat = "YYYY-MM-DD"
